$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the daily conversion message in A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 4.24 = 16441.9 pesos`n✅ 16441.9 pesos = 4.23 = 959.59 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$wsHoja1.Range("A1").Value = $newText

# --- Sheet "tasas": update rate figures ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 235.8
$wsTasas.Range("N12").Value = 3889.5
$wsTasas.Range("O12").Value = 227
